# "Top 100.xlsx" refresh: per-artist aggregate stats (plays + audio-feature
# averages in columns M:T) were recalculated after the backing "plays" data
# was updated and all helper computations were folded into BaseClasses.
# This workbook has no formulas (every cell is a cached literal), so the
# refreshed numbers are written directly onto the affected cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 0.3676923076923078
$ws.Range("P4").Value = 0.05545454545454546
$ws.Range("N5").Value = 0.2073913043478261
$ws.Range("P5").Value = 0.04565217391304349
$ws.Range("M6").Value = 202
$ws.Range("N7").Value = 0.2241666666666666
$ws.Range("O7").Value = 0.4050000000000001
$ws.Range("M8").Value = 31
$ws.Range("O8").Value = 0.44
$ws.Range("M10").Value = 29
$ws.Range("P11").Value = 0.05222222222222222
$ws.Range("O13").Value = 0.4100000000000001
$ws.Range("P13").Value = 0.03333333333333335
$ws.Range("R14").Value = 0.8014285714285715
$ws.Range("S14").Value = 0.2064285714285714
$ws.Range("P15").Value = 0.05950000000000002
$ws.Range("O17").Value = 0.5946666666666666
$ws.Range("R19").Value = 0.5861538461538461
$ws.Range("P22").Value = 0.0953846153846154
$ws.Range("Q22").Value = 0.4392307692307692
$ws.Range("R22").Value = 0.8461538461538461
$ws.Range("T22").Value = 0.2553846153846154
$ws.Range("M24").Value = 127
$ws.Range("R24").Value = 0.6093333333333334
$ws.Range("O25").Value = 0.4540909090909092
$ws.Range("P25").Value = 0.05318181818181819
$ws.Range("Q25").Value = 0.4390909090909092
$ws.Range("P26").Value = 0.03599999999999999
$ws.Range("R26").Value = 0.75
$ws.Range("P28").Value = 0.1776923076923077
$ws.Range("Q28").Value = 0.4015384615384616
$ws.Range("R28").Value = 0.6653846153846155
$ws.Range("O29").Value = 0.5433333333333333
$ws.Range("P29").Value = 0.07433333333333336
$ws.Range("R30").Value = 0.4485714285714286
$ws.Range("Q31").Value = 0.3391666666666666
$ws.Range("R31").Value = 0.5966666666666667
$ws.Range("Q32").Value = 0.49
$ws.Range("R32").Value = 0.8088235294117645
$ws.Range("S32").Value = 0.1870588235294117
$ws.Range("R33").Value = 0.4745454545454546
$ws.Range("M34").Value = 15
$ws.Range("P34").Value = 0.05692307692307694
$ws.Range("R34").Value = 0.7061538461538461
$ws.Range("O35").Value = 0.5713333333333332
$ws.Range("M36").Value = 9
$ws.Range("N36").Value = 0.4218749999999999
$ws.Range("Q36").Value = 0.4012500000000001
$ws.Range("N37").Value = 0.4318750000000001
$ws.Range("Q37").Value = 0.3593750000000001
$ws.Range("Q38").Value = 0.4081818181818182
$ws.Range("R38").Value = 0.8236363636363637
$ws.Range("P39").Value = 0.4185714285714285
$ws.Range("R39").Value = 0.5952380952380952
$ws.Range("P40").Value = 0.03857142857142859
$ws.Range("N41").Value = 0.07916666666666668
$ws.Range("R42").Value = 0.9662500000000001
$ws.Range("N43").Value = 0.2845454545454545
$ws.Range("P44").Value = 0.03300000000000001
$ws.Range("Q45").Value = 0.1638461538461539
$ws.Range("R45").Value = 0.4584615384615385
$ws.Range("M46").Value = 68
$ws.Range("M47").Value = 23
$ws.Range("Q47").Value = 0.6927272727272727
$ws.Range("R47").Value = 0.7836363636363636
$ws.Range("S48").Value = 0.1566666666666667
$ws.Range("O49").Value = 0.5141666666666667
$ws.Range("Q49").Value = 0.4125000000000001
$ws.Range("P50").Value = 0.0675
$ws.Range("T50").Value = 0.5075000000000001
$ws.Range("Q51").Value = 0.3176470588235295
$ws.Range("R52").Value = 0.9225
$ws.Range("M53").Value = 79
$ws.Range("O53").Value = 0.6671428571428571
$ws.Range("P53").Value = 0.08785714285714286
$ws.Range("R55").Value = 0.5666666666666667
$ws.Range("M56").Value = 6
$ws.Range("M58").Value = 7
$ws.Range("N58").Value = 0.03333333333333334
$ws.Range("O58").Value = 0.4233333333333333
$ws.Range("P58").Value = 0.04777777777777778
$ws.Range("Q58").Value = 0.1533333333333333
$ws.Range("R58").Value = 0.7255555555555556
$ws.Range("S58").Value = 0.1644444444444445
$ws.Range("T58").Value = 0.27
$ws.Range("M59").Value = 10
$ws.Range("O59").Value = 0.4991666666666667
$ws.Range("R59").Value = 0.8275000000000001
$ws.Range("N61").Value = 0.0742857142857143
$ws.Range("O61").Value = 0.5735714285714286
$ws.Range("P61").Value = 0.05928571428571429
$ws.Range("Q61").Value = 0.5785714285714286
$ws.Range("R61").Value = 0.8757142857142857
$ws.Range("M63").Value = 63
$ws.Range("P63").Value = 0.09833333333333333
$ws.Range("Q64").Value = 0.4407692307692308
$ws.Range("P68").Value = 0.03272727272727272
$ws.Range("O69").Value = 0.654
$ws.Range("Q69").Value = 0.57
$ws.Range("T71").Value = 0.2421428571428572
$ws.Range("O72").Value = 0.4687499999999999
$ws.Range("P72").Value = 0.04875000000000001
$ws.Range("Q72").Value = 0.405
$ws.Range("P73").Value = 0.0468421052631579
$ws.Range("R73").Value = 0.5910526315789475
